# ddc algorithm downloaded from moodle
# Adds new Work-plan log entries (rows 16-19), updates the "LECTURE:" note,
# fixes up the related styling/row-height, and nudges the remembered
# selections on the "Work plan" and "Facts" sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Work plan")
$ws2 = $wb.Worksheets.Item("Facts")

# --- New log text (order matters: it drives shared-string table layout) ---
$ws1.Range("E17").Value = "Learn to measure execution time in matlab"
$ws1.Range("E16").Value = "LECTURE: code planning"
$ws1.Range("F16").Value = "I stayed till the end without break. Find out that I can use clutering to define which model is closiest to average(mean)"
$ws1.Range("E18").Value = "Map the matrix indexes to coordinates when drawing the plot"
$ws1.Range("E19").Value = "Save ensemble to simple mean file"

# --- Formatting for the new cells (reuse existing styles, don't invent new ones) ---
# B17/B18/C18/D18 -> same "date" style as B16/C16
$ws1.Range("B16").Copy()
$ws1.Range("B17").PasteSpecial(-4122)
$ws1.Range("B18").PasteSpecial(-4122)
$ws1.Range("C18").PasteSpecial(-4122)
$ws1.Range("D18").PasteSpecial(-4122)

# F16 -> same wrap-text style as E15
$ws1.Range("E15").Copy()
$ws1.Range("F16").PasteSpecial(-4122)

$ws1.Application.CutCopyMode = $false

# --- Date values for the newly added date cells ---
$ws1.Range("B17").Value2 = 43892
$ws1.Range("B18").Value2 = 43892
$ws1.Range("C18").Value2 = 43892
$ws1.Range("D18").Value2 = 43892

# --- Row 16 grows to fit the new wrapped F16 note ---
$ws1.Rows("16").RowHeight = 43.2

# --- Selection bookkeeping: remembered active cell moves with the edits ---
$ws1.Range("F19").Select()
$ws2.Range("B16").Select()
$ws1.Activate()
